$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New computed values for cells A2:A86 (row -> value), per the commit's updated results
$newValues = @(
    5.763841670213395,
    9.126103885849659,
    10.36343549267141,
    11.40455399623508,
    6.035988099228689,
    11.68693180418234,
    8.589748610734119,
    7.436663384409371,
    7.815448069575837,
    7.945335257886313,
    2.451116197982941,
    6.178423552399238,
    4.518903128704466,
    2.100334423147075,
    2.343551991077874,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197,
    5.064110380611197
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $newValues[$i]
}
